$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire rows for "RM 232" (row 26) and "SC 92" (row 28).
# Delete the lower row first so the earlier row index stays valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(26).Delete()

# Fill-in / clear individual cell values (after row deletion, rows below
# row 26 have shifted up by one; rows below the original row 28 have
# shifted up by two).

# C3 (row 3, "RM 8"): empty -> 11.2
$ws.Range("C3").Value = 11.2

# D4 (row 4, "RM 9"): -15.4 -> empty
$ws.Range("D4").Value = ""

# C5 (row 5, "RM 14"): 12.3 -> empty
$ws.Range("C5").Value = ""

# D9 (row 9, "RM 42"): empty -> -14.5
$ws.Range("D9").Value = -14.5

# D10 (row 10, "RM 52 a"): empty -> -14.7
$ws.Range("D10").Value = -14.7

# D13 (row 13, "RM 88"): -13.9 -> empty
$ws.Range("D13").Value = ""

# D14 (row 14, "RM 90"): -13.1 -> empty
$ws.Range("D14").Value = ""

# C21 (row 21, "RM 135"): empty -> 12.7
$ws.Range("C21").Value = 12.7

# C23 (row 23, "RM 140"): 12.2 -> empty
$ws.Range("C23").Value = ""

# C32 (row 32 after deletions, "SC 193"): empty -> 10.5
$ws.Range("C32").Value = 10.5
